$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8562.375
$ws.Range("I62").Value = 4166.6665
$ws.Range("K62").Value = 4166.6665
$ws.Range("M62").Value = -3542.6665
$ws.Range("H65").Value = 8562.375
$ws.Range("I65").Value = 4166.6665
$ws.Range("K65").Value = 20833.3325
$ws.Range("M65").Value = -17713.3325
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 4637
$ws.Range("I113").Value = 6002.5
$ws.Range("J113").Value = 1906
$ws.Range("K113").Value = 6002.5
$ws.Range("L113").Value = 1906
$ws.Range("M113").Value = -2748.5
$ws.Range("N113").Value = -8414
$ws.Range("H141").Value = 1982.1666
$ws.Range("I141").Value = 1883.3334
$ws.Range("J141").Value = 2081
$ws.Range("K141").Value = 5650.0002
$ws.Range("L141").Value = 6243
$ws.Range("M141").Value = -470.0002000000004
$ws.Range("N141").Value = -16603

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2685.6
$ws.Range("I2").Value = 5677.5
$ws.Range("J2").Value = 1937.625
$ws.Range("K2").Value = 5677.5
$ws.Range("L2").Value = 1937.625
$ws.Range("M2").Value = -5564.5
$ws.Range("N2").Value = -2163.625
$ws.Range("H110").Value = 125002280
$ws.Range("I110").Value = 166668620
$ws.Range("K110").Value = 166668620
$ws.Range("M110").Value = -166666575
$ws.Range("H116").Value = 2685.6
$ws.Range("I116").Value = 5677.5
$ws.Range("J116").Value = 1937.625
$ws.Range("K116").Value = 5677.5
$ws.Range("L116").Value = 1937.625
$ws.Range("M116").Value = -3383.5
$ws.Range("N116").Value = -6525.625

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2685.6
$ws.Range("I3").Value = 5677.5
$ws.Range("J3").Value = 1937.625
$ws.Range("K3").Value = 5677.5
$ws.Range("L3").Value = 1937.625
$ws.Range("M3").Value = -5563.5
$ws.Range("N3").Value = -2165.625
$ws.Range("H82").Value = 21960.666
$ws.Range("I82").Value = 12219.75
$ws.Range("J82").Value = 99888
$ws.Range("K82").Value = 12219.75
$ws.Range("L82").Value = 99888
$ws.Range("M82").Value = -11836.75
$ws.Range("N82").Value = -100654
$ws.Range("H85").Value = 21960.666
$ws.Range("I85").Value = 12219.75
$ws.Range("J85").Value = 99888
$ws.Range("K85").Value = 12219.75
$ws.Range("L85").Value = 99888
$ws.Range("M85").Value = -10893.75
$ws.Range("N85").Value = -102540

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2285.5
$ws.Range("I16").Value = 2107
$ws.Range("J16").Value = 2999.5
$ws.Range("K16").Value = 2107
$ws.Range("L16").Value = 2999.5
$ws.Range("M16").Value = -1820
$ws.Range("N16").Value = -3573.5
$ws.Range("H31").Value = 4686.5
$ws.Range("J31").Value = 6249.85
$ws.Range("L31").Value = 6249.85
$ws.Range("N31").Value = -6839.85
$ws.Range("H34").Value = 4686.5
$ws.Range("J34").Value = 6249.85
$ws.Range("L34").Value = 6249.85
$ws.Range("N34").Value = -6653.85
$ws.Range("H51").Value = 98765
$ws.Range("J51").Value = 98765
$ws.Range("L51").Value = 98765
$ws.Range("N51").Value = -100237
$ws.Range("H58").Value = 3132.0625
$ws.Range("I58").Value = 1760.75
$ws.Range("J58").Value = 7246
$ws.Range("K58").Value = 1760.75
$ws.Range("L58").Value = 7246
$ws.Range("M58").Value = -1557.75
$ws.Range("N58").Value = -7652
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 20752.143
$ws.Range("I60").Value = 7750
$ws.Range("K60").Value = 7750
$ws.Range("M60").Value = -7239
$ws.Range("H61").Value = 98765
$ws.Range("J61").Value = 98765
$ws.Range("L61").Value = 98765
$ws.Range("N61").Value = -99461
$ws.Range("H68").Value = 74190.75
$ws.Range("J68").Value = 74190.75
$ws.Range("L68").Value = 74190.75
$ws.Range("N68").Value = -75688.75
$ws.Range("H71").Value = 74190.75
$ws.Range("J71").Value = 74190.75
$ws.Range("L71").Value = 222572.25
$ws.Range("N71").Value = -230060.25
$ws.Range("H92").Value = 54999.5
$ws.Range("J92").Value = 54999.5
$ws.Range("L92").Value = 54999.5
$ws.Range("N92").Value = -59991.5
$ws.Range("H113").Value = 2285.5
$ws.Range("I113").Value = 2107
$ws.Range("J113").Value = 2999.5
$ws.Range("K113").Value = 2107
$ws.Range("L113").Value = 2999.5
$ws.Range("M113").Value = 63
$ws.Range("N113").Value = -7339.5
$ws.Range("H136").Value = 3132.0625
$ws.Range("I136").Value = 1760.75
$ws.Range("J136").Value = 7246
$ws.Range("K136").Value = 5282.25
$ws.Range("L136").Value = 21738
$ws.Range("M136").Value = -2732.25
$ws.Range("N136").Value = -26838

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1858
$ws.Range("I5").Value = 1726.8572
$ws.Range("K5").Value = 5180.571599999999
$ws.Range("M5").Value = -5068.571599999999
$ws.Range("H113").Value = 1276.9
$ws.Range("I113").Value = 665
$ws.Range("J113").Value = 1429.875
$ws.Range("K113").Value = 1995
$ws.Range("L113").Value = 4289.625
$ws.Range("M113").Value = 175
$ws.Range("N113").Value = -8629.625
$ws.Range("H135").Value = 1858
$ws.Range("I135").Value = 1726.8572
$ws.Range("K135").Value = 15541.7148
$ws.Range("M135").Value = -13006.7148

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 60097.5
$ws.Range("J57").Value = 60097.5
$ws.Range("L57").Value = 60097.5
$ws.Range("N57").Value = -61737.5
$ws.Range("H92").Value = 9999.166999999999
$ws.Range("J92").Value = 9999.166999999999
$ws.Range("L92").Value = 9999.166999999999
$ws.Range("N92").Value = -13743.167
$ws.Range("H99").Value = 2447.4285
$ws.Range("I99").Value = 2447.4285
$ws.Range("K99").Value = 2447.4285
$ws.Range("M99").Value = -201.4285
$ws.Range("H101").Value = 45995
$ws.Range("J101").Value = 45995
$ws.Range("L101").Value = 45995
$ws.Range("N101").Value = -52485
$ws.Range("H113").Value = 4199
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4535
$ws.Range("I40").Value = 3380
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 3380
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -3244
$ws.Range("N40").Value = -8272
$ws.Range("H46").Value = 7127.722
$ws.Range("I46").Value = 5299.909
$ws.Range("K46").Value = 5299.909
$ws.Range("M46").Value = -5111.909
$ws.Range("H55").Value = 1198.7142
$ws.Range("I55").Value = 1551.8
$ws.Range("J55").Value = 1002.55554
$ws.Range("K55").Value = 1551.8
$ws.Range("L55").Value = 1002.55554
$ws.Range("M55").Value = -1378.8
$ws.Range("N55").Value = -1348.55554
$ws.Range("H61").Value = 166671570
$ws.Range("I61").Value = 200004880
$ws.Range("K61").Value = 200004880
$ws.Range("M61").Value = -200004678
$ws.Range("H99").Value = 65000
$ws.Range("I99").Value = 65000
$ws.Range("K99").Value = 65000
$ws.Range("M99").Value = -62005
$ws.Range("H113").Value = 166671570
$ws.Range("I113").Value = 200004880
$ws.Range("K113").Value = 200004880
$ws.Range("M113").Value = -200002710

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3760.625
$ws.Range("J81").Value = 7668
$ws.Range("L81").Value = 15336
$ws.Range("N81").Value = -17458
$ws.Range("H84").Value = 3760.625
$ws.Range("J84").Value = 7668
$ws.Range("L84").Value = 76680
$ws.Range("N84").Value = -87288
$ws.Range("H107").Value = 333333340
$ws.Range("I107").Value = 333333340
$ws.Range("K107").Value = 1000000020
$ws.Range("M107").Value = -999998100
